$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("O2").Value = 30.31
$ws1.Range("O12").Value = "1 de 10"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F2").Value = 30.31
$ws2.Range("F12").Value = 140.37

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D2").Value = 1769.08
$ws3.Range("E2").Value = -1769.08

$ws3.Range("D3").Value = 911.6799999999999
$ws3.Range("E3").Value = 12811.66
$ws3.Range("F3").Value = 0.06643280717376382

$ws3.Range("D4").Value = 2680.76
$ws3.Range("E4").Value = 11042.58
$ws3.Range("F4").Value = 0.1953431161801719
